$d = $word.ActiveDocument

# --- Portuguese "Programa" paragraph: split single run of text into 5 numbered
# items separated by manual line breaks (w:br) ---
$d.Content.Find.Execute(
    "diferentes processos2. Processos",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "diferentes processos^l2. Processos",
    2) | Out-Null

$d.Content.Find.Execute(
    "alterações de alimentos 3. Processos",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "alterações de alimentos ^l3. Processos",
    2) | Out-Null

$d.Content.Find.Execute(
    "produtos desidratados 4. Discussão",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "produtos desidratados ^l4. Discussão",
    2) | Out-Null

$d.Content.Find.Execute(
    "importância industrial5. Bioenergia",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "importância industrial^l5. Bioenergia",
    2) | Out-Null

# --- English (italic) "Programa" paragraph ---
$d.Content.Find.Execute(
    "different processes2. Biochemical",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "different processes^l2. Biochemical",
    2) | Out-Null

$d.Content.Find.Execute(
    "food conservation/modifications3. Biochemical",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "food conservation/modifications^l3. Biochemical",
    2) | Out-Null

$d.Content.Find.Execute(
    "dehydrated products4. Discussion",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "dehydrated products^l4. Discussion",
    2) | Out-Null

$d.Content.Find.Execute(
    "industrial importance5. Bioenergy",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "industrial importance^l5. Bioenergy",
    2) | Out-Null

# --- Bibliografia paragraph: split the 3 references onto separate lines ---
$d.Content.Find.Execute(
    "9788521313823.LIMA, U. A.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "9788521313823.^lLIMA, U. A.",
    2) | Out-Null

$d.Content.Find.Execute(
    "9788521214571.Moraes, I. O.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "9788521214571.^lMoraes, I. O.",
    2) | Out-Null
